{"js": "// Replace the date line and each two-digit\u00f7one-digit division problem cell\n// with its updated value. Each old string is unique within the document, so a\n// plain exact-text search (matchCase, no wildcards) safely targets exactly\n// one run each.\nconst replacements = [\n  [\n    \"2026-01-04 Sunday\",\n    \"2026-01-05 Monday\"\n  ],\n  [\n    \"24\u00f77=3, 3\",\n    \"90\u00f74=22, 2\"\n  ],\n  [\n    \"48\u00f73=16, 0\",\n    \"79\u00f72=39, 1\"\n  ],\n  [\n    \"97\u00f79=10, 7\",\n    \"19\u00f77=2, 5\"\n  ],\n  [\n    \"96\u00f76=16, 0\",\n    \"75\u00f79=8, 3\"\n  ],\n  [\n    \"34\u00f76=5, 4\",\n    \"98\u00f73=32, 2\"\n  ],\n  [\n    \"28\u00f79=3, 1\",\n    \"39\u00f72=19, 1\"\n  ],\n  [\n    \"21\u00f76=3, 3\",\n    \"41\u00f76=6, 5\"\n  ],\n  [\n    \"40\u00f74=10, 0\",\n    \"41\u00f72=20, 1\"\n  ],\n  [\n    \"77\u00f76=12, 5\",\n    \"68\u00f75=13, 3\"\n  ],\n  [\n    \"92\u00f74=23, 0\",\n    \"79\u00f77=11, 2\"\n  ],\n  [\n    \"48\u00f76=8, 0\",\n    \"82\u00f78=10, 2\"\n  ],\n  [\n    \"67\u00f73=22, 1\",\n    \"23\u00f79=2, 5\"\n  ],\n  [\n    \"80\u00f76=13, 2\",\n    \"74\u00f76=12, 2\"\n  ],\n  [\n    \"71\u00f73=23, 2\",\n    \"96\u00f72=48, 0\"\n  ],\n  [\n    \"71\u00f75=14, 1\",\n    \"34\u00f73=11, 1\"\n  ],\n  [\n    \"29\u00f72=14, 1\",\n    \"13\u00f74=3, 1\"\n  ],\n  [\n    \"18\u00f72=9, 0\",\n    \"56\u00f79=6, 2\"\n  ],\n  [\n    \"57\u00f76=9, 3\",\n    \"13\u00f76=2, 1\"\n  ],\n  [\n    \"90\u00f75=18, 0\",\n    \"19\u00f73=6, 1\"\n  ],\n  [\n    \"64\u00f77=9, 1\",\n    \"87\u00f77=12, 3\"\n  ],\n  [\n    \"44\u00f75=8, 4\",\n    \"80\u00f78=10, 0\"\n  ],\n  [\n    \"38\u00f72=19, 0\",\n    \"91\u00f74=22, 3\"\n  ],\n  [\n    \"47\u00f77=6, 5\",\n    \"43\u00f79=4, 7\"\n  ],\n  [\n    \"26\u00f75=5, 1\",\n    \"91\u00f72=45, 1\"\n  ],\n  [\n    \"59\u00f72=29, 1\",\n    \"16\u00f78=2, 0\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date line and each two-digit division problem cell\n# to its new value. Each old string is unique in the document, so Find/Replace\n# (exact text, not whole document formatting) safely retargets exactly one run.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2026-01-04 Sunday\", \"2026-01-05 Monday\"),\n  @(\"24\u00f77=3, 3\", \"90\u00f74=22, 2\"),\n  @(\"48\u00f73=16, 0\", \"79\u00f72=39, 1\"),\n  @(\"97\u00f79=10, 7\", \"19\u00f77=2, 5\"),\n  @(\"96\u00f76=16, 0\", \"75\u00f79=8, 3\"),\n  @(\"34\u00f76=5, 4\", \"98\u00f73=32, 2\"),\n  @(\"28\u00f79=3, 1\", \"39\u00f72=19, 1\"),\n  @(\"21\u00f76=3, 3\", \"41\u00f76=6, 5\"),\n  @(\"40\u00f74=10, 0\", \"41\u00f72=20, 1\"),\n  @(\"77\u00f76=12, 5\", \"68\u00f75=13, 3\"),\n  @(\"92\u00f74=23, 0\", \"79\u00f77=11, 2\"),\n  @(\"48\u00f76=8, 0\", \"82\u00f78=10, 2\"),\n  @(\"67\u00f73=22, 1\", \"23\u00f79=2, 5\"),\n  @(\"80\u00f76=13, 2\", \"74\u00f76=12, 2\"),\n  @(\"71\u00f73=23, 2\", \"96\u00f72=48, 0\"),\n  @(\"71\u00f75=14, 1\", \"34\u00f73=11, 1\"),\n  @(\"29\u00f72=14, 1\", \"13\u00f74=3, 1\"),\n  @(\"18\u00f72=9, 0\", \"56\u00f79=6, 2\"),\n  @(\"57\u00f76=9, 3\", \"13\u00f76=2, 1\"),\n  @(\"90\u00f75=18, 0\", \"19\u00f73=6, 1\"),\n  @(\"64\u00f77=9, 1\", \"87\u00f77=12, 3\"),\n  @(\"44\u00f75=8, 4\", \"80\u00f78=10, 0\"),\n  @(\"38\u00f72=19, 0\", \"91\u00f74=22, 3\"),\n  @(\"47\u00f77=6, 5\", \"43\u00f79=4, 7\"),\n  @(\"26\u00f75=5, 1\", \"91\u00f72=45, 1\"),\n  @(\"59\u00f72=29, 1\", \"16\u00f78=2, 0\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}\n\n"}
